{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the edit (see commit message \"Update ho\u00e0n ch\u1ec9nh slide + b\u00e1o\n// c\u00e1o ch\u01b0a c\u00f3 demo\"):\n//   1. Extend the bullet \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1\n//      n\u00e0o \u0111\u00f3.\" with an additional clause, making it:\n//      \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1 n\u00e0o \u0111\u00f3, con s\u1ed1 n\u00e0y\n//      ph\u1ea3i c\u00f3 \u00fd ngh\u0129a nh\u1ea5t \u0111\u1ecbnh trong ch\u01b0\u01a1ng tr\u00ecnh c\u1ee7a ch\u00fang ta.\"\n//   2. Add a brand-new bullet right after it (same list/level):\n//      \"N\u1ebfu m\u1ed9t l\u1edbp l\u00e0 Singleton, ch\u00fang ta c\u00f3 th\u1ec3 g\u1eb7p kh\u00f3 kh\u0103n khi t\u1ea1o l\u1edbp\n//      con c\u1ee7a l\u1edbp \u0111\u00f3.\"\n//   3. The document's \"_GoBack\" bookmark (Word's \"last edit location\"\n//      marker) follows the new text, i.e. it moves from the end of the\n//      document to right after \"...nh\u1ea5t \u0111\u1ecbnh \" (where typing stopped).\n\nconst body = context.document.body;\n\nconst originalSentence =\n  \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1 n\u00e0o \u0111\u00f3.\";\nconst updatedSentence =\n  \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1 n\u00e0o \u0111\u00f3, con s\u1ed1 n\u00e0y ph\u1ea3i c\u00f3 \u00fd ngh\u0129a nh\u1ea5t \u0111\u1ecbnh trong ch\u01b0\u01a1ng tr\u00ecnh c\u1ee7a ch\u00fang ta.\";\nconst newBulletText =\n  \"N\u1ebfu m\u1ed9t l\u1edbp l\u00e0 Singleton, ch\u00fang ta c\u00f3 th\u1ec3 g\u1eb7p kh\u00f3 kh\u0103n khi t\u1ea1o l\u1edbp con c\u1ee7a l\u1edbp \u0111\u00f3.\";\n\n// 1) Find the target paragraph/run via search.\nconst results = body.search(originalSentence, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target sentence to update.\");\n}\n\nconst targetRange = results.items[0];\nconst targetParagraph = targetRange.paragraphs.getFirst();\n\n// 2) Replace the sentence text in place (keeps the run's formatting).\ntargetRange.insertText(updatedSentence, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Insert the new bullet paragraph right after, inheriting the same\n//    list numbering/formatting as the paragraph it follows.\nconst newParagraph = targetParagraph.insertParagraph(\n  newBulletText,\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 4) Move the \"_GoBack\" bookmark to sit right after \"nh\u1ea5t \u0111\u1ecbnh \" (i.e.\n//    right before \"trong ch\u01b0\u01a1ng tr\u00ecnh c\u1ee7a ch\u00fang ta.\"), matching where\n//    Word leaves it after the edit. Remove it from its old location (end\n//    of document) first, if present.\nconst goBackRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\ngoBackRange.load(\"isNullObject\");\nawait context.sync();\n\nif (!goBackRange.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nconst markerResults = body.search(\"nh\u1ea5t \u0111\u1ecbnh \", { matchCase: true });\nmarkerResults.load(\"items\");\nawait context.sync();\n\nif (markerResults.items.length > 0) {\n  const markerRange = markerResults.items[0];\n  const insertionPoint = markerRange.getRange(Word.RangeLocation.after);\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Summary of the edit (see commit message \"Update ho\u00e0n ch\u1ec9nh slide + b\u00e1o\n# c\u00e1o ch\u01b0a c\u00f3 demo\"):\n#   1. Extend the bullet \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1\n#      n\u00e0o \u0111\u00f3.\" with an additional clause, making it:\n#      \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1 n\u00e0o \u0111\u00f3, con s\u1ed1 n\u00e0y\n#      ph\u1ea3i c\u00f3 \u00fd ngh\u0129a nh\u1ea5t \u0111\u1ecbnh trong ch\u01b0\u01a1ng tr\u00ecnh c\u1ee7a ch\u00fang ta.\"\n#   2. Add a brand-new bullet right after it (same list/level):\n#      \"N\u1ebfu m\u1ed9t l\u1edbp l\u00e0 Singleton, ch\u00fang ta c\u00f3 th\u1ec3 g\u1eb7p kh\u00f3 kh\u0103n khi t\u1ea1o l\u1edbp\n#      con c\u1ee7a l\u1edbp \u0111\u00f3.\"\n#   3. The document's \"_GoBack\" bookmark (Word's \"last edit location\"\n#      marker) follows the new text, i.e. it moves from the end of the\n#      document to right after \"...nh\u1ea5t \u0111\u1ecbnh \" (where typing stopped).\n\n$d = $word.ActiveDocument\n\n$originalSentence = \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1 n\u00e0o \u0111\u00f3.\"\n$updatedSentence = \"Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng th\u1ec3 hi\u1ec7n c\u1ee7a l\u1edbp \u1edf m\u1ed9t con s\u1ed1 n\u00e0o \u0111\u00f3, con s\u1ed1 n\u00e0y ph\u1ea3i c\u00f3 \u00fd ngh\u0129a nh\u1ea5t \u0111\u1ecbnh trong ch\u01b0\u01a1ng tr\u00ecnh c\u1ee7a ch\u00fang ta.\"\n$newBulletText = \"N\u1ebfu m\u1ed9t l\u1edbp l\u00e0 Singleton, ch\u00fang ta c\u00f3 th\u1ec3 g\u1eb7p kh\u00f3 kh\u0103n khi t\u1ea1o l\u1edbp con c\u1ee7a l\u1edbp \u0111\u00f3.\"\n$markerText = \"nh\u1ea5t \u0111\u1ecbnh \"\n\n# 1) Find and replace the sentence text (keeps the run's formatting) and\n#    open up a new paragraph right after it.\n$rng = $d.Content\n$found = $rng.Find.Execute($originalSentence)\nif ($found) {\n  $rng.Text = $updatedSentence\n  $rng.InsertParagraphAfter()\n}\n\n# 2) Fill in the newly-opened paragraph with the new bullet text. It\n#    already inherited the same list numbering/formatting from the\n#    paragraph above (InsertParagraphAfter duplicates pPr).\n$paras = $d.Paragraphs\n$count = $paras.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $paras.Item($i)\n  if ($p.Range.Text -like \"*Gi\u1edbi h\u1ea1n s\u1ed1 l\u01b0\u1ee3ng*\") {\n    $newPara = $p.Next()\n    $newPara.Range.InsertBefore($newBulletText)\n    break\n  }\n}\n\n# 3) Move the \"_GoBack\" bookmark to sit right after \"nh\u1ea5t \u0111\u1ecbnh \" (i.e.\n#    right before \"trong ch\u01b0\u01a1ng tr\u00ecnh c\u1ee7a ch\u00fang ta.\"), matching where\n#    Word leaves it after the edit. Remove it from its old location (end\n#    of document) first, if present.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$markerRange = $d.Content\n$foundMarker = $markerRange.Find.Execute($markerText)\nif ($foundMarker) {\n  $markerRange.Collapse(0)\n  $d.Bookmarks.Add(\"_GoBack\", $markerRange)\n}\n"}
